$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase5")
$ws.Activate()

# Update B8:B41 (rows for gestweek_conception 6..39) from 0.02 to 0.01
for ($r = 8; $r -le 41; $r++) {
    $ws.Cells.Item($r, 2).Value = 0.01
}

# Update the selected cell/range on the sheet to E34
$ws.Range("E34").Select()
